# CONTRAT_EMBAUCHE_A_ESSAI.docx edit
#
# 1) The "Direction :" field placeholder "${direction}" becomes
#    "${direction_acceuil}", split (as real Word does when text is typed
#    mid-run) into three runs: "${direction", "_acceuil" and "}".
# 2) The lone "_GoBack" bookmark (previously sitting right after the
#    "Objet :" colon, near the top of the document) is Word's
#    last-edit-position marker, so it moves to sit right where the user's
#    cursor ended up: between the newly typed "_acceuil" and the closing
#    "}".
#
# Word only ever keeps a single "_GoBack" bookmark, so re-adding it at the
# new location automatically removes it from its old spot (handled below
# with Bookmarks.Add, which repositions the bookmark with id 0 rather
# than creating a second one).

$d = $word.ActiveDocument

# Locate the exact "${direction}" placeholder (there are similarly named
# "${direction_sc}" / "${direction_acceuil}" placeholders elsewhere in the
# document, so search for the full, closed token to land on the right one).
$rng = $d.Content
$found = $rng.Find.Execute("`${direction}")
if (-not $found) {
    throw "Could not find the `${direction} placeholder"
}

$openStart = $rng.Start          # start of "${direction}"
$splitPos1 = $openStart + 11     # "${direction" is 11 characters -> right before "}"
$splitPos2 = $splitPos1 + 8      # "_acceuil" is 8 characters -> right before "}"

# Type "_acceuil" right before the closing brace, turning "${direction}"
# into "${direction_acceuil}".
$insertRange = $d.Range($splitPos1, $splitPos1)
$insertRange.InsertAfter("_acceuil")

# Drop a throwaway bookmark exactly at the "${direction" / "_acceuil"
# boundary purely so the two stay separate runs (mirroring the two
# distinctly-typed runs seen in the authored edit); it is removed again
# immediately below.
$markerRange = $d.Range($splitPos1, $splitPos1)
$d.Bookmarks.Add("zzz_tmp_split", $markerRange)

# Place (move) the real "_GoBack" bookmark at the final cursor position,
# right after "_acceuil" and before the closing "}".
$goBackRange = $d.Range($splitPos2, $splitPos2)
$d.Bookmarks.Add("_GoBack", $goBackRange)

# Drop the scratch bookmark now that the run split it was protecting has
# already been committed.
$d.Bookmarks("zzz_tmp_split").Delete()
